{"js": "// Replace each math-expression / date cell's text with its updated value.\n// old -> new pairs correspond 1:1 and in document order to the cells in the body.\nconst replacements = [\n  [\"2024-08-11 Sunday\", \"2024-08-12 Monday\"],\n  [\"24+26=\", \"54-51=\"],\n  [\"55+5=\", \"19+71=\"],\n  [\"47+7=\", \"72+15=\"],\n  [\"14+45=\", \"8+91=\"],\n  [\"79-37=\", \"56+15=\"],\n  [\"91+3=\", \"99-72=\"],\n  [\"23-8=\", \"48-37=\"],\n  [\"56-53=\", \"12+20=\"],\n  [\"26+66=\", \"6+68=\"],\n  [\"66-49=\", \"79-0=\"],\n  [\"2+33=\", \"93-66=\"],\n  [\"7+64=\", \"78-76=\"],\n  [\"82-67=\", \"17+30=\"],\n  [\"74+23=\", \"65-14=\"],\n  [\"29-6=\", \"24+34=\"],\n  [\"4+55=\", \"24-12=\"],\n  [\"44+40=\", \"55-43=\"],\n  [\"36-21=\", \"49+1=\"],\n  [\"69-52=\", \"27-2=\"],\n  [\"64-20=\", \"65-28=\"],\n  [\"56-38=\", \"40+36=\"],\n  [\"72+11=\", \"87-43=\"],\n  [\"52-10=\", \"1+25=\"],\n  [\"95-34=\", \"5+60=\"],\n  [\"63+11=\", \"69+30=\"],\n  [\"43+46=\", \"30+19=\"],\n  [\"22+53=\", \"99-62=\"],\n  [\"45+16=\", \"31+62=\"],\n  [\"66-13=\", \"92-73=\"],\n  [\"29+13=\", \"95-72=\"],\n  [\"54+22=\", \"0+41=\"],\n  [\"1+18=\", \"12+30=\"],\n  [\"61+31=\", \"63-45=\"],\n  [\"8+68=\", \"26+39=\"],\n  [\"49+13=\", \"26-18=\"],\n  [\"53-34=\", \"23+53=\"],\n  [\"86+0=\", \"99-96=\"],\n  [\"28+30=\", \"60-40=\"],\n  [\"95-2=\", \"49-43=\"],\n  [\"58+22=\", \"97-68=\"],\n  [\"5+65=\", \"85-53=\"],\n  [\"56+23=\", \"0+49=\"],\n  [\"13-3=\", \"76-40=\"],\n  [\"96-25=\", \"35+49=\"],\n  [\"54+41=\", \"45+41=\"],\n  [\"0+62=\", \"86-77=\"],\n  [\"88-47=\", \"72-38=\"],\n  [\"73-26=\", \"5+2=\"],\n  [\"67-4=\", \"69-26=\"],\n  [\"92-81=\", \"18+72=\"],\n  [\"34+12=\", \"52+5=\"],\n  [\"29+67=\", \"4+77=\"],\n  [\"53-29=\", \"3+91=\"],\n  [\"35+29=\", \"45+23=\"],\n  [\"32+61=\", \"71-55=\"],\n  [\"57+37=\", \"81-58=\"],\n  [\"58-44=\", \"80-36=\"],\n  [\"9+76=\", \"70-46=\"],\n  [\"64+17=\", \"95-39=\"],\n  [\"74-67=\", \"27-27=\"],\n  [\"83+5=\", \"3+77=\"],\n  [\"14+44=\", \"41+1=\"],\n  [\"39+5=\", \"58+32=\"],\n  [\"61-61=\", \"85-79=\"],\n  [\"85-32=\", \"95-37=\"],\n  [\"68-46=\", \"88-54=\"],\n  [\"80-37=\", \"54-8=\"],\n  [\"35-19=\", \"28+69=\"],\n  [\"26+46=\", \"70-7=\"],\n  [\"80-29=\", \"29+29=\"],\n  [\"25+40=\", \"65+19=\"],\n  [\"89-67=\", \"45+35=\"],\n  [\"62-60=\", \"12+57=\"],\n  [\"36-24=\", \"23+61=\"],\n  [\"91-49=\", \"15-5=\"],\n  [\"73-44=\", \"25-18=\"],\n  [\"29-15=\", \"20+75=\"],\n  [\"74-36=\", \"30-29=\"],\n  [\"29-21=\", \"84-63=\"],\n  [\"64-8=\", \"96-21=\"],\n  [\"33+26=\", \"5+31=\"],\n  [\"23+33=\", \"23+45=\"],\n  [\"54-45=\", \"10+22=\"],\n  [\"24+16=\", \"18+28=\"],\n  [\"92-92=\", \"91-26=\"],\n  [\"66+3=\", \"75-33=\"],\n  [\"73-38=\", \"69+21=\"],\n  [\"23-2=\", \"5+21=\"],\n  [\"93-58=\", \"52+41=\"],\n  [\"69-58=\", \"64-29=\"],\n  [\"44-18=\", \"51-47=\"],\n  [\"62+28=\", \"64+5=\"],\n  [\"96-89=\", \"83-61=\"],\n  [\"48+41=\", \"46-23=\"],\n  [\"73-39=\", \"27+68=\"],\n  [\"96-68=\", \"76-68=\"],\n  [\"57-50=\", \"38+47=\"],\n  [\"47+6=\", \"94-92=\"],\n  [\"47-13=\", \"71+21=\"],\n  [\"23+29=\", \"10+82=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "# Replace each math-expression / date cell's text with its updated value.\n# old -> new pairs correspond 1:1 and in document order to the cells in the body.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-08-11 Sunday\", \"2024-08-12 Monday\"),\n  @(\"24+26=\", \"54-51=\"),\n  @(\"55+5=\", \"19+71=\"),\n  @(\"47+7=\", \"72+15=\"),\n  @(\"14+45=\", \"8+91=\"),\n  @(\"79-37=\", \"56+15=\"),\n  @(\"91+3=\", \"99-72=\"),\n  @(\"23-8=\", \"48-37=\"),\n  @(\"56-53=\", \"12+20=\"),\n  @(\"26+66=\", \"6+68=\"),\n  @(\"66-49=\", \"79-0=\"),\n  @(\"2+33=\", \"93-66=\"),\n  @(\"7+64=\", \"78-76=\"),\n  @(\"82-67=\", \"17+30=\"),\n  @(\"74+23=\", \"65-14=\"),\n  @(\"29-6=\", \"24+34=\"),\n  @(\"4+55=\", \"24-12=\"),\n  @(\"44+40=\", \"55-43=\"),\n  @(\"36-21=\", \"49+1=\"),\n  @(\"69-52=\", \"27-2=\"),\n  @(\"64-20=\", \"65-28=\"),\n  @(\"56-38=\", \"40+36=\"),\n  @(\"72+11=\", \"87-43=\"),\n  @(\"52-10=\", \"1+25=\"),\n  @(\"95-34=\", \"5+60=\"),\n  @(\"63+11=\", \"69+30=\"),\n  @(\"43+46=\", \"30+19=\"),\n  @(\"22+53=\", \"99-62=\"),\n  @(\"45+16=\", \"31+62=\"),\n  @(\"66-13=\", \"92-73=\"),\n  @(\"29+13=\", \"95-72=\"),\n  @(\"54+22=\", \"0+41=\"),\n  @(\"1+18=\", \"12+30=\"),\n  @(\"61+31=\", \"63-45=\"),\n  @(\"8+68=\", \"26+39=\"),\n  @(\"49+13=\", \"26-18=\"),\n  @(\"53-34=\", \"23+53=\"),\n  @(\"86+0=\", \"99-96=\"),\n  @(\"28+30=\", \"60-40=\"),\n  @(\"95-2=\", \"49-43=\"),\n  @(\"58+22=\", \"97-68=\"),\n  @(\"5+65=\", \"85-53=\"),\n  @(\"56+23=\", \"0+49=\"),\n  @(\"13-3=\", \"76-40=\"),\n  @(\"96-25=\", \"35+49=\"),\n  @(\"54+41=\", \"45+41=\"),\n  @(\"0+62=\", \"86-77=\"),\n  @(\"88-47=\", \"72-38=\"),\n  @(\"73-26=\", \"5+2=\"),\n  @(\"67-4=\", \"69-26=\"),\n  @(\"92-81=\", \"18+72=\"),\n  @(\"34+12=\", \"52+5=\"),\n  @(\"29+67=\", \"4+77=\"),\n  @(\"53-29=\", \"3+91=\"),\n  @(\"35+29=\", \"45+23=\"),\n  @(\"32+61=\", \"71-55=\"),\n  @(\"57+37=\", \"81-58=\"),\n  @(\"58-44=\", \"80-36=\"),\n  @(\"9+76=\", \"70-46=\"),\n  @(\"64+17=\", \"95-39=\"),\n  @(\"74-67=\", \"27-27=\"),\n  @(\"83+5=\", \"3+77=\"),\n  @(\"14+44=\", \"41+1=\"),\n  @(\"39+5=\", \"58+32=\"),\n  @(\"61-61=\", \"85-79=\"),\n  @(\"85-32=\", \"95-37=\"),\n  @(\"68-46=\", \"88-54=\"),\n  @(\"80-37=\", \"54-8=\"),\n  @(\"35-19=\", \"28+69=\"),\n  @(\"26+46=\", \"70-7=\"),\n  @(\"80-29=\", \"29+29=\"),\n  @(\"25+40=\", \"65+19=\"),\n  @(\"89-67=\", \"45+35=\"),\n  @(\"62-60=\", \"12+57=\"),\n  @(\"36-24=\", \"23+61=\"),\n  @(\"91-49=\", \"15-5=\"),\n  @(\"73-44=\", \"25-18=\"),\n  @(\"29-15=\", \"20+75=\"),\n  @(\"74-36=\", \"30-29=\"),\n  @(\"29-21=\", \"84-63=\"),\n  @(\"64-8=\", \"96-21=\"),\n  @(\"33+26=\", \"5+31=\"),\n  @(\"23+33=\", \"23+45=\"),\n  @(\"54-45=\", \"10+22=\"),\n  @(\"24+16=\", \"18+28=\"),\n  @(\"92-92=\", \"91-26=\"),\n  @(\"66+3=\", \"75-33=\"),\n  @(\"73-38=\", \"69+21=\"),\n  @(\"23-2=\", \"5+21=\"),\n  @(\"93-58=\", \"52+41=\"),\n  @(\"69-58=\", \"64-29=\"),\n  @(\"44-18=\", \"51-47=\"),\n  @(\"62+28=\", \"64+5=\"),\n  @(\"96-89=\", \"83-61=\"),\n  @(\"48+41=\", \"46-23=\"),\n  @(\"73-39=\", \"27+68=\"),\n  @(\"96-68=\", \"76-68=\"),\n  @(\"57-50=\", \"38+47=\"),\n  @(\"47+6=\", \"94-92=\"),\n  @(\"47-13=\", \"71+21=\"),\n  @(\"23+29=\", \"10+82=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}"}
